# Clean-up pass: remove the three screenshot paragraphs, the three
# "────" divider lines, and the nine empty spacer paragraphs
# (pPr/spacing w:before="40") that sat right after each code table.
#
# We scan Paragraphs back-to-front (after collecting the indices we want
# to remove) so that deleting one paragraph never invalidates the index
# of another one still pending deletion.

$d = $word.ActiveDocument

$toDelete = New-Object System.Collections.ArrayList

$total = $d.Paragraphs.Count
for ($i = 1; $i -le $total; $i++) {
    $p = $d.Paragraphs($i)
    $rng = $p.Range

    # Tables contribute a zero-length "boundary" paragraph at their end
    # in the Paragraphs collection; it isn't real document content, skip it.
    if ($rng.Start -eq $rng.End) {
        continue
    }

    $text = $rng.Text
    $trimmed = $text.Trim([char]13, [char]7)

    $hasDrawing = $rng.InlineShapes.Count -gt 0
    $isSeparator = ($trimmed.Length -gt 0) -and ($trimmed.Substring(0, 1) -eq [char]0x2500)
    $isEmptySpacer = ($trimmed.Length -eq 0) -and (-not $hasDrawing) -and ($p.Format.SpaceBefore -eq 2)

    if ($hasDrawing -or $isSeparator -or $isEmptySpacer) {
        [void]$toDelete.Add($i)
    }
}

for ($j = $toDelete.Count - 1; $j -ge 0; $j--) {
    $idx = $toDelete[$j]
    $d.Paragraphs($idx).Range.Delete()
}
